$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new week's price record was added for "Cebollín" at
# "Terminal Hortofrutícola Agro Chillán". It is inserted as the new row 4,
# pushing the previous rows 4-10 down to rows 5-11 (the sheet keeps its
# existing newest-last-of-group ordering pattern).
$ws.Rows.Item(4).Insert()

# Populate the newly inserted row 4 with the new weekly record.
$ws.Range("A4").Value = 7
$ws.Range("B4").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C4").Value = "Ñuble"
$ws.Range("D4").Value = 44662
$ws.Range("E4").Value = 16
$ws.Range("F4").Value = 100112037
$ws.Range("G4").Value = "Cebollín"
$ws.Range("H4").Value = "Sin especificar"
$ws.Range("I4").Value = "Primera"
$ws.Range("J4").Value = 200
$ws.Range("K4").Value = 8000
$ws.Range("L4").Value = 8500
$ws.Range("M4").Value = 8250
$ws.Range("N4").Value = "$/paquete 36 unidades"
$ws.Range("O4").Value = "Región Metropolitana"
$ws.Range("P4").Value = 229
$ws.Range("Q4").Value = 36
$ws.Range("R4").Value = "Hortaliza"
